$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = (Get-Date -Year 2021 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 2900
$ws.Range("O2").Value = 3000
$ws.Range("P2").Value = 2950
$ws.Range("S2").Value = 2950
$ws.Range("D3").Value = (Get-Date -Year 2021 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N3").Value = 2600
$ws.Range("O3").Value = 2600
$ws.Range("P3").Value = 2600
$ws.Range("S3").Value = 2600
$ws.Range("D4").Value = (Get-Date -Year 2021 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M4").Value = 150
$ws.Range("P4").Value = 25467
$ws.Range("S4").Value = 2547
$ws.Range("D5").Value = (Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 29000
$ws.Range("P5").Value = 28500
$ws.Range("S5").Value = 2850
$ws.Range("D6").Value = (Get-Date -Year 2021 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N6").Value = 29000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29500
$ws.Range("Q6").Value = '$/bandeja 10 kilos'
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 2950
$ws.Range("T6").Value = 10
$ws.Range("D7").Value = (Get-Date -Year 2021 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("S7").Value = 2950
$ws.Range("D8").Value = (Get-Date -Year 2021 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N8").Value = 1900
$ws.Range("O8").Value = 2000
$ws.Range("P8").Value = 1950
$ws.Range("S8").Value = 1950
$ws.Range("D9").Value = (Get-Date -Year 2021 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 1700
$ws.Range("O9").Value = 1700
$ws.Range("P9").Value = 1700
$ws.Range("Q9").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S9").Value = 1700
$ws.Range("T9").Value = 1
$ws.Range("D10").Value = (Get-Date -Year 2021 -Month 10 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 25000
$ws.Range("O10").Value = 26000
$ws.Range("P10").Value = 25600
$ws.Range("S10").Value = 2560
$ws.Range("D11").Value = (Get-Date -Year 2021 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N11").Value = 2000
$ws.Range("O11").Value = 2100
$ws.Range("P11").Value = 2050
$ws.Range("Q11").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S11").Value = 2050
$ws.Range("T11").Value = 1
$ws.Range("D12").Value = (Get-Date -Year 2022 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N12").Value = 29000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 29500
$ws.Range("Q12").Value = '$/bandeja 10 kilos'
$ws.Range("S12").Value = 2950
$ws.Range("T12").Value = 10
$ws.Range("L13").Value = 'Primera'
$ws.Range("N13").Value = 2200
$ws.Range("O13").Value = 2200
$ws.Range("P13").Value = 2200
$ws.Range("S13").Value = 2200
$ws.Range("D14").Value = (Get-Date -Year 2021 -Month 11 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L14").Value = 'Segunda'
$ws.Range("N14").Value = 1800
$ws.Range("O14").Value = 1800
$ws.Range("P14").Value = 1800
$ws.Range("S14").Value = 1800
$ws.Range("D15").Value = (Get-Date -Year 2022 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L15").Value = 'Especial'
$ws.Range("N15").Value = 22000
$ws.Range("O15").Value = 22000
$ws.Range("P15").Value = 22000
$ws.Range("Q15").Value = '$/bandeja 10 kilos'
$ws.Range("S15").Value = 2200
$ws.Range("T15").Value = 10
$ws.Range("D16").Value = (Get-Date -Year 2022 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 19000
$ws.Range("P16").Value = 19000
$ws.Range("Q16").Value = '$/bandeja 10 kilos'
$ws.Range("S16").Value = 1900
$ws.Range("T16").Value = 10
$ws.Range("D17").Value = (Get-Date -Year 2022 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 17000
$ws.Range("P17").Value = 17000
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("S17").Value = 1700
$ws.Range("T17").Value = 10
$ws.Range("D18").Value = (Get-Date -Year 2022 -Month 9 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M18").Value = 80
$ws.Range("O18").Value = 26000
$ws.Range("P18").Value = 25500
$ws.Range("S18").Value = 2550
$ws.Range("D19").Value = (Get-Date -Year 2021 -Month 10 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 2600
$ws.Range("O19").Value = 2600
$ws.Range("P19").Value = 2600
$ws.Range("Q19").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S19").Value = 2600
$ws.Range("T19").Value = 1
$ws.Range("D20").Value = (Get-Date -Year 2021 -Month 10 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = 2400
$ws.Range("O20").Value = 2400
$ws.Range("P20").Value = 2400
$ws.Range("Q20").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S20").Value = 2400
$ws.Range("T20").Value = 1
$ws.Range("D21").Value = (Get-Date -Year 2021 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 26000
$ws.Range("O21").Value = 27000
$ws.Range("P21").Value = 26500
$ws.Range("S21").Value = 2650
$ws.Range("D22").Value = (Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 2700
$ws.Range("O22").Value = 2800
$ws.Range("P22").Value = 2750
$ws.Range("Q22").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S22").Value = 2750
$ws.Range("T22").Value = 1
$ws.Range("D23").Value = (Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L23").Value = 'Segunda'
$ws.Range("N23").Value = 2500
$ws.Range("O23").Value = 2500
$ws.Range("P23").Value = 2500
$ws.Range("Q23").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S23").Value = 2500
$ws.Range("T23").Value = 1
$ws.Range("D24").Value = (Get-Date -Year 2020 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 100
$ws.Range("O24").Value = 18000
$ws.Range("P24").Value = 17500
$ws.Range("Q24").Value = '$/bandeja 8 kilos'
$ws.Range("S24").Value = 2188
$ws.Range("T24").Value = 8
$ws.Range("D25").Value = (Get-Date -Year 2020 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 15000
$ws.Range("O25").Value = 15000
$ws.Range("P25").Value = 15000
$ws.Range("Q25").Value = '$/bandeja 8 kilos'
$ws.Range("S25").Value = 1875
$ws.Range("T25").Value = 8
$ws.Range("D26").Value = (Get-Date -Year 2021 -Month 11 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 80
$ws.Range("N26").Value = 25000
$ws.Range("O26").Value = 26000
$ws.Range("P26").Value = 25375
$ws.Range("Q26").Value = '$/bandeja 10 kilos'
$ws.Range("S26").Value = 2538
$ws.Range("T26").Value = 10
$ws.Range("D27").Value = (Get-Date -Year 2021 -Month 10 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 26000
$ws.Range("P27").Value = 25500
$ws.Range("S27").Value = 2550
$ws.Range("D28").Value = (Get-Date -Year 2021 -Month 9 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N28").Value = 3200
$ws.Range("O28").Value = 3300
$ws.Range("P28").Value = 3250
$ws.Range("Q28").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R28").Value = 'Provincia del Elquí'
$ws.Range("S28").Value = 3250
$ws.Range("T28").Value = 1
$ws.Range("D29").Value = (Get-Date -Year 2021 -Month 9 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N29").Value = 30000
$ws.Range("O29").Value = 31000
$ws.Range("P29").Value = 30500
$ws.Range("Q29").Value = '$/bandeja 10 kilos'
$ws.Range("S29").Value = 3050
$ws.Range("T29").Value = 10
$ws.Range("D30").Value = (Get-Date -Year 2021 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 25000
$ws.Range("O30").Value = 27000
$ws.Range("P30").Value = 26000
$ws.Range("Q30").Value = '$/bandeja 10 kilos'
$ws.Range("S30").Value = 2600
$ws.Range("T30").Value = 10
